$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A8').Value = 'Volume 30   Number  35'
$ws.Range('C9').Value = 'Report Covering the Week  8/28/2023  Through  9/3/2023'
$ws.Range('D14').Value = 1
$ws.Range('D14').NumberFormat = '#,##0'
$ws.Range('E14').Value = -100
$ws.Range('E14').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('G14').Value = 1
$ws.Range('G14').NumberFormat = '#,##0'
$ws.Range('H14').Value = -100
$ws.Range('H14').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('J14').Value = 5
$ws.Range('K14').Value = -80
$ws.Range('N14').Value = -88.888888888888
$ws.Range('M15').Value = 55.555555555555
$ws.Range('N15').Value = -50
$ws.Range('C16').Value = 2
$ws.Range('D16').Value = 6
$ws.Range('E16').Value = -66.666666666666
$ws.Range('F16').Value = 15
$ws.Range('G16').Value = 14
$ws.Range('H16').Value = 7.142857142857
$ws.Range('I16').Value = 130
$ws.Range('J16').Value = 171
$ws.Range('K16').Value = -23.976608187134
$ws.Range('L16').Value = 31.313131313131
$ws.Range('M16').Value = 5.691056910569
$ws.Range('N16').Value = -79.166666666666
$ws.Range('C17').Value = 7
$ws.Range('E17').Value = 16.666666666666
$ws.Range('G17').Value = 21
$ws.Range('H17').Value = -38.095238095238
$ws.Range('I17').Value = 162
$ws.Range('J17').Value = 137
$ws.Range('K17').Value = 18.248175182481
$ws.Range('L17').Value = 23.664122137404
$ws.Range('M17').Value = 54.285714285714
$ws.Range('N17').Value = -59.5
$ws.Range('C18').Value = 5
$ws.Range('D18').Value = 13
$ws.Range('E18').Value = -61.538461538461
$ws.Range('F18').Value = 20
$ws.Range('G18').Value = 32
$ws.Range('H18').Value = -37.5
$ws.Range('I18').Value = 183
$ws.Range('J18').Value = 264
$ws.Range('K18').Value = -30.681818181818
$ws.Range('L18').Value = -5.670103092783
$ws.Range('M18').Value = 17.307692307692
$ws.Range('N18').Value = -68.983050847457
$ws.Range('C19').Value = 20
$ws.Range('D19').Value = 24
$ws.Range('E19').Value = -16.666666666666
$ws.Range('F19').Value = 82
$ws.Range('G19').Value = 99
$ws.Range('H19').Value = -17.171717171717
$ws.Range('I19').Value = 681
$ws.Range('J19').Value = 704
$ws.Range('K19').Value = -3.267045454545
$ws.Range('L19').Value = 49.342105263157
$ws.Range('M19').Value = 30.210325047801
$ws.Range('N19').Value = -34.582132564841
$ws.Range('D20').Value = 1
$ws.Range('D20').NumberFormat = '#,##0'
$ws.Range('E20').Value = 0
$ws.Range('E20').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('I20').Value = 31
$ws.Range('J20').Value = 30
$ws.Range('K20').Value = 3.333333333333
$ws.Range('L20').Value = -55.714285714285
$ws.Range('M20').Value = -8.823529411764
$ws.Range('N20').Value = -91.621621621621
$ws.Range('C21').Value = 35
$ws.Range('D21').Value = 51
$ws.Range('E21').Value = -31.372549019607
$ws.Range('F21').Value = 135
$ws.Range('G21').Value = 170
$ws.Range('H21').Value = -20.588235294117
$ws.Range('I21').Value = 1202
$ws.Range('J21').Value = 1324
$ws.Range('K21').Value = -9.214501510574
$ws.Range('L21').Value = 24.173553719008
$ws.Range('M21').Value = 26.128016789087
$ws.Range('N21').Value = -60.74461136512
$ws.Range('C23').Value = 2
$ws.Range('D23').Value = 1
$ws.Range('E23').Value = 100
$ws.Range('F23').Value = 10
$ws.Range('G23').Value = 7
$ws.Range('H23').Value = 42.857142857142
$ws.Range('I23').Value = 93
$ws.Range('J23').Value = 94
$ws.Range('K23').Value = -1.063829787234
$ws.Range('L23').Value = -31.111111111111
$ws.Range('M23').Value = 8.13953488372
$ws.Range('C24').Value = 36
$ws.Range('E24').Value = -10
$ws.Range('G24').Value = 164
$ws.Range('H24').Value = -31.097560975609
$ws.Range('I24').Value = 986
$ws.Range('J24').Value = 1520
$ws.Range('K24').Value = -35.131578947368
$ws.Range('L24').Value = 39.660056657223
$ws.Range('M24').Value = -14.926660914581
$ws.Range('C25').Value = 11
$ws.Range('D25').Value = 6
$ws.Range('E25').Value = 83.333333333333
$ws.Range('F25').Value = 42
$ws.Range('G25').Value = 41
$ws.Range('H25').Value = 2.439024390243
$ws.Range('I25').Value = 325
$ws.Range('J25').Value = 332
$ws.Range('K25').Value = -2.108433734939
$ws.Range('L25').Value = 31.048387096774
$ws.Range('M25').Value = 0.308641975308
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = '0'
$ws.Range('C27').NumberFormat = 'General'
$ws.Range('D27').Value = 1
$ws.Range('E27').Value = -100
$ws.Range('F27').Value = 2
$ws.Range('G27').Value = 6
$ws.Range('H27').Value = -66.666666666666
$ws.Range('J27').Value = 66
$ws.Range('K27').Value = -46.969696969697
$ws.Range('L27').Value = -7.894736842105
$ws.Range('D28').Value = 1
$ws.Range('D28').NumberFormat = '#,##0'
$ws.Range('E28').Value = -100
$ws.Range('E28').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('G28').Value = 1
$ws.Range('G28').NumberFormat = '#,##0'
$ws.Range('H28').Value = -100
$ws.Range('H28').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('J28').Value = 4
$ws.Range('K28').Value = -25
$ws.Range('D29').Value = 1
$ws.Range('D29').NumberFormat = '#,##0'
$ws.Range('E29').Value = -100
$ws.Range('E29').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('G29').Value = 1
$ws.Range('G29').NumberFormat = '#,##0'
$ws.Range('H29').Value = -100
$ws.Range('H29').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('J29').Value = 4
$ws.Range('K29').Value = -25
